$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 (Transistors (x4) / BC337): update the Comments column (E)
$ws.Range("E9").Value = "Use Damo's to test"

# Row 11 (Motors (x4) / H107-A03): update the From column (D)
$ws.Range("D11").Value = "eBay"
